# The workbook's single sheet "Avverkningsanmälningar" has a "Förändrad"
# (changed/modified) date column C for every data row (rows 2-436). The
# commit updates that date by one day (serial 46061 -> 46062, i.e.
# 2026-02-08 -> 2026-02-09) for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("C2:C436")
$range.Value2 = 46062
